$d = $word.ActiveDocument

# Locate the exact text "Hands-on Practice (" that needs to be restructured.
$find = $d.Content
$found = $find.Find.Execute("Hands-on Practice (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $find.Start

    $part1 = "Hands-on Practice"
    $part2 = " on W3schools and TUF Platform"
    $part3 = " ("

    # Replace the original "Hands-on Practice (" text with the full new combined
    # text. This necessarily merges with whatever run(s) immediately follow
    # (e.g. the existing "8" / " hours)" runs), so afterwards we re-split the
    # touched region back into the desired individual runs.
    $whole = $d.Range($start, $find.End)
    $whole.Text = $part1 + $part2 + $part3

    $p1End = $start + $part1.Length
    $p2End = $p1End + $part2.Length
    $p3End = $p2End + $part3.Length

    # Toggling a character formatting property on/off forces the engine to
    # split the run boundary at that point without leaving any residual
    # formatting difference behind.
    $split1 = $d.Range($start, $p1End)
    $split1.Bold = 1
    $split1.Bold = 0

    $split2 = $d.Range($p1End, $p2End)
    $split2.Bold = 1
    $split2.Bold = 0

    $split3 = $d.Range($p2End, $p3End)
    $split3.Bold = 1
    $split3.Bold = 0

    # The text right after $p3End used to be the pre-existing "8" run
    # followed by " hours)"; our edit above merged "8" into the same run as
    # " (" (since they shared identical formatting). Split them back apart
    # so "8" and " hours)" remain their own separate, untouched-looking runs.
    # Search only a small local window right after $p3End so we can't
    # accidentally latch onto an unrelated digit elsewhere in the document.
    $digitRange = $d.Range($p3End, $p3End + 10)
    $digitFound = $digitRange.Find.Execute("8", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($digitFound -and $digitRange.Start -eq $p3End) {
        $digitEnd = $digitRange.End
        $split4 = $d.Range($p3End, $digitEnd)
        $split4.Bold = 1
        $split4.Bold = 0
    }
}
